# Insert a new data row above the current row 32 (this pushes the
# existing rows 32..114 down to 33..115, preserving all of their data
# and formatting), then populate the newly inserted row 32 with the
# new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 32, shifting rows 32:114 down to 33:115.
$ws.Rows.Item(32).Insert()

# Fill in the new row 32 with the new record.
$ws.Range("A32").Value2 = 11
$ws.Range("B32").Value2 = 'Vega Monumental Concepción'
$ws.Range("C32").Value2 = 'Bíobío'
$ws.Range("D32").Value2 = 44659
$ws.Range("E32").Value2 = 8
$ws.Range("F32").Value2 = 100112032
$ws.Range("G32").Value2 = 'Zapallo italiano'
$ws.Range("H32").Value2 = 'Sin especificar'
$ws.Range("I32").Value2 = 'Primera'
$ws.Range("J32").Value2 = 350
$ws.Range("K32").Value2 = 8500
$ws.Range("L32").Value2 = 9000
$ws.Range("M32").Value2 = 8714
$ws.Range("N32").Value2 = '$/caja 50 unidades'
$ws.Range("O32").Value2 = 'Región Metropolitana'
$ws.Range("P32").Value2 = 174
$ws.Range("Q32").Value2 = 50
$ws.Range("R32").Value2 = 'Hortaliza'

# Match the date-cell number formatting used by the rest of column D.
$ws.Range("D32").NumberFormat = $ws.Range("D33").NumberFormat
